# Generate Report for Handoff
#
# - Flip the localization status from "In Translation" to "Ready for
#   handoff" everywhere it is shown (Overview sheet per-language status
#   cells, and the per-language "Status" column on the zh-cn / de-de
#   detail sheets).
# - Refresh the associated handoff timestamps to the new generation time.
# - Re-size the datetime columns so the longer "Ready for handoff" /
#   refreshed-timestamp text still fits (mirrors an Excel AutoFit).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# --- Refreshed handoff / HO-xliff-generate timestamps ---
$zhcn.Range("H2").Value = "2016-08-29 00:57:25"
$dede.Range("H2").Value = "2016-08-29 00:57:29"
$overview.Range("G2").Value = "2016-08-29 00:57:29"

# --- Widen the datetime columns to fit the new text (AutoFit-style) ---
# The COM ColumnWidth setter here quantizes to 1/6-character increments
# (pixel granularity), so feed it the character-width input whose
# quantized result lands closest to the desired rendered width.
$newWidth = 98 / 6
$overview.Range("E1").ColumnWidth = $newWidth
$overview.Range("F1").ColumnWidth = $newWidth
$zhcn.Range("C1").ColumnWidth = $newWidth
$dede.Range("C1").ColumnWidth = $newWidth
